$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column (H) — copy the header style used by the other
# header cells (bold, centered, bordered) from G1, then set the text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data values for the Save column
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
